$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 334; this shifts existing rows 334:417 down to 335:418
# and extends the used range / dimension to A1:R418, matching the target diff.
$ws.Rows(334).Insert()

# Populate the newly inserted row 334 with the new data record
$ws.Range("A334").Value = 8
$ws.Range("B334").Value = "Terminal La Palmera de La Serena"
$ws.Range("C334").Value = "Coquimbo"
$ws.Range("D334").Value = 44855
$ws.Range("E334").Value = 4
$ws.Range("F334").Value = 100112032
$ws.Range("G334").Value = "Zapallo italiano"
$ws.Range("H334").Value = "Sin especificar"
$ws.Range("I334").Value = "Primera"
$ws.Range("J334").Value = 600
$ws.Range("K334").Value = 11500
$ws.Range("L334").Value = 12000
$ws.Range("M334").Value = 11750
$ws.Range("N334").Value = "`$/caja 50 unidades"
$ws.Range("O334").Value = "Región de Arica y Parinacota"
$ws.Range("P334").Value = 235
$ws.Range("Q334").Value = 50
$ws.Range("R334").Value = "Hortaliza"
